$d = $word.ActiveDocument
$sec = $d.Sections(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
$targets = @(
    @{ HF = $sec.Footers(1); NewName = "image2.png" },
    @{ HF = $sec.Footers(2); NewName = "image2.png" },
    @{ HF = $sec.Headers(1); NewName = "image1.jpg" },
    @{ HF = $sec.Headers(2); NewName = "image1.jpg" }
)

foreach ($t in $targets) {
    $hf = $t.HF
    if ($hf.Exists) {
        $ishp = $hf.Range.InlineShapes(1)
        # InlineShape has no settable Name property in the Word object model;
        # renaming requires the round-trip through a floating Shape.
        $shp = $ishp.ConvertToShape()
        $shp.Name = $t.NewName
        [void]$shp.ConvertToInlineShape()
    }
}

Write-Host "Renamed inline pictures in headers/footers."
